# Applies the "Atualizacao de bases das ligas" update to the Denmark Superligaen sheet.
# The update swaps the full data (B,F..AC columns) between certain pairs/cycles of
# rows (newer odds-feed rows were re-ordered), and applies a couple of standalone
# odds corrections on two additional rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: row -> row whose B..AC content it should receive (i.e. a permutation).
# Expressed as disjoint cycles found in the source data.
$cycles = New-Object System.Collections.ArrayList
[void]$cycles.Add(@(3,4))
[void]$cycles.Add(@(9,10))
[void]$cycles.Add(@(15,16))
[void]$cycles.Add(@(26,28,31))
[void]$cycles.Add(@(27,30,29))
[void]$cycles.Add(@(51,52))
[void]$cycles.Add(@(57,58))
[void]$cycles.Add(@(75,76))
[void]$cycles.Add(@(81,82))
[void]$cycles.Add(@(86,87,88))
[void]$cycles.Add(@(89,90))
[void]$cycles.Add(@(142,143))
[void]$cycles.Add(@(172,173))
[void]$cycles.Add(@(202,203))

# Columns B (2) through AC (29) hold: id, Div, Div Original Name, Date... wait -
# actually B..AC here = external id (B) .. PL_AhUnder (AC); keep A (row idx),
# C, D, E untouched.
$firstCol = 2   # B
$lastCol  = 29  # AC

foreach ($cycle in $cycles) {
    # Snapshot the B:AC row block for every row in this cycle BEFORE any writes.
    $snapshot = @{}
    foreach ($r in $cycle) {
        $rng = $ws.Range($ws.Cells.Item($r, $firstCol), $ws.Cells.Item($r, $lastCol))
        $snapshot[$r] = $rng.Value()
    }

    # Rotate: row[i] receives the content that was in row[i+1] (cyclically),
    # matching the source-data permutation (e.g. cycle (26,28,31) means
    # 26<-28, 28<-31, 31<-26).
    $n = $cycle.Count
    for ($i = 0; $i -lt $n; $i++) {
        $destRow = $cycle[$i]
        $srcRow  = $cycle[($i + 1) % $n]
        $rng = $ws.Range($ws.Cells.Item($destRow, $firstCol), $ws.Cells.Item($destRow, $lastCol))
        $rng.Value = $snapshot[$srcRow]
    }
}

# Standalone odds corrections (no row swap involved).
$ws.Range("R204").Value = 1.88
$ws.Range("S204").Value = 2.02
$ws.Range("U204").Value = 1.875
$ws.Range("V204").Value = 1.975

$ws.Range("R205").Value = 1.86
$ws.Range("S205").Value = 2.04
$ws.Range("U205").Value = 2.025
$ws.Range("V205").Value = 1.825
